$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each destination row (2-20), the D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) values are
# replaced with the values that originally belonged to another row (a full
# permutation of the 19 data rows). Snapshot the original values first so the
# permutation can be applied without clobbering source data before it is read.

$rows = 2..20

# Snapshot original values for columns D, J, K, L, M, P (columns 4, 10, 11, 12, 13, 16)
$origD = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origP = @{}

foreach ($r in $rows) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value2
    $origJ[$r] = $ws.Cells.Item($r, 10).Value2
    $origK[$r] = $ws.Cells.Item($r, 11).Value2
    $origL[$r] = $ws.Cells.Item($r, 12).Value2
    $origM[$r] = $ws.Cells.Item($r, 13).Value2
    $origP[$r] = $ws.Cells.Item($r, 16).Value2
}

# destination row -> source row (the row whose original values now populate it)
$permutation = @{
    2  = 5
    3  = 10
    4  = 17
    5  = 9
    6  = 11
    7  = 2
    8  = 16
    9  = 3
    10 = 15
    11 = 6
    12 = 7
    13 = 19
    14 = 13
    15 = 18
    16 = 4
    17 = 14
    18 = 8
    19 = 20
    20 = 12
}

foreach ($r in $rows) {
    $src = $permutation[$r]
    $ws.Cells.Item($r, 4).Value = $origD[$src]
    $ws.Cells.Item($r, 10).Value = $origJ[$src]
    $ws.Cells.Item($r, 11).Value = $origK[$src]
    $ws.Cells.Item($r, 12).Value = $origL[$src]
    $ws.Cells.Item($r, 13).Value = $origM[$src]
    $ws.Cells.Item($r, 16).Value = $origP[$src]
}
